# Trade #45 closed at 2026-02-17 13:28:01 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# ---- Summary sheet ----
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1197.51   # Current Capital
$summary.Range("B4").Value = -2.48     # Total P&L $
$summary.Range("B6").Value = 45        # Total Trades
$summary.Range("B8").Value = 24        # Losing Trades
$summary.Range("B9").Value = 40        # Win Rate %

# ---- Strategy Status sheet (MarketMaking row) ----
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 97.51000000000001  # Capital
$status.Range("D4").Value = 45                 # Trades
$status.Range("E4").Value = -2.48              # P&L $
$status.Range("F4").Value = -2.49              # P&L %
$status.Range("G4").Value = 40                 # Win Rate %

# ---- Append the new closed trade (Trade #45) to both the ----
# "All Trades" log and the strategy-specific "MarketMaking" log.
$sheetNames = @("All Trades", "MarketMaking")
foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("A46").Value = 45

    # Date column looks like "2026-02-17" which Excel would otherwise
    # auto-convert to a date serial; force it to stay plain text to match
    # the rest of the column, then strip the temporary format again.
    $dateCell = $ws.Range("B46")
    $dateCell.NumberFormat = "@"
    $dateCell.Value = "2026-02-17"
    $dateCell.Style = "Normal"

    $ws.Range("C46").Value = "13:27:54"
    $ws.Range("D46").Value = "MarketMaking"
    $ws.Range("E46").Value = "UP"
    $ws.Range("F46").Value = 0.883866
    $ws.Range("G46").Value = 0.83
    $ws.Range("H46").Value = "CLOSED"
    $ws.Range("I46").Value = -6.0944
    $ws.Range("J46").Value = -0.05
    $ws.Range("K46").Value = 97.51000000000001
    $ws.Range("L46").Value = 0
    $ws.Range("M46").Value = 0
    $ws.Range("N46").Value = 0.6
    $ws.Range("O46").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P46").Value = "early_exit"
    $ws.Range("Q46").Value = 0.13
}
